$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the column headers (row 1) - "Modification des noms de colonnes"
$ws.Range("A1").Value = "identifiant_indic"   # was: indic_id
# B1 "zone_id" is unchanged
$ws.Range("C1").Value = "date_valeur"         # was: metric_date
$ws.Range("D1").Value = "type_valeur"         # was: metric_type
$ws.Range("E1").Value = "valeur"              # was: metric_value

# Move / persist the selected cell to E2
[void]$ws.Range("E2").Select()

# Keep the 1904 date system flag explicit (matches the semantics of the
# workbookPr date1904 attribute, even though its value itself is unchanged).
$wb.Date1904 = $false
